$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9
$ws.Cells.Item($row, 1).Value = "'13"
$ws.Cells.Item($row, 2).Value = "newrelic_alert_channel"
$ws.Cells.Item($row, 3).Value = "open"
$ws.Cells.Item($row, 4).Value = "2025-03-24T09:05:44Z"
$ws.Cells.Item($row, 5).Value = "bug"
